$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.726.88'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '2.619.61'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("D9").Value = '2.618.98'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +9.71%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("D16").Value = '3.097.68'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").Value = '67.699.63'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '2.621.34'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '365.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.92%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("E27").Value = '  +2.73%  '
$ws.Range("D28").Value = '2.745.46'
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '584.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("E34").Value = '  -1.88%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.45'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '155.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.57%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("E48").Value = '  -6.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.624'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.08%  '
